# Auto-generated Excel COM-interop edit script.
# Applies the cryptos.xlsx price / 1h-volume refresh described by the
# commit diff ("Updated cryptos list ... with GitHub Actions"):
#   - Column D ("Price") and column E ("Volume(1h)") values are refreshed
#     for most rows.
#   - Rows 38-41 are additionally re-ranked: Kaspa/Dai/Fetch.AI swap rows
#     (Fetch.AI moves to row 38, VeChain stays at 39, Kaspa moves to row 40,
#     Dai moves to row 41), each bringing its own Link/Price/Volume along.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D stores prices as plain text (e.g. "61.00", "0.0000346",
# "68.382.87") even though many of them look numeric. Writing a plain
# numeric-looking string into Range.Value lets Excel reinterpret it as a
# number (dropping trailing zeros / switching to scientific notation), so
# we briefly force Text format on the price column, write the values, then
# restore the "Normal" style so the cells end up with no style override -
# matching the original sheet exactly - while the stored text is preserved.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

# --- Column D ("Price") updates ------------------------------------------
$ws.Range("D2").Value = '68.382.87'
$ws.Range("D3").Value = '3.916.66'
$ws.Range("D5").Value = '485.68'
$ws.Range("D6").Value = '146.01'
$ws.Range("D7").Value = '0.623'
$ws.Range("D8").Value = '0.997'
$ws.Range("D9").Value = '0.734'
$ws.Range("D11").Value = '0.0000346'
$ws.Range("D12").Value = '43.17'
$ws.Range("D14").Value = '4.538.66'
$ws.Range("D15").Value = '3.877.63'
$ws.Range("D16").Value = '14.34'
$ws.Range("D18").Value = '19.98'
$ws.Range("D20").Value = '68.414.48'
$ws.Range("D21").Value = '432.82'
$ws.Range("D22").Value = '15.20'
$ws.Range("D24").Value = '88.16'
$ws.Range("D25").Value = '11.51'
$ws.Range("D28").Value = '37.95'
$ws.Range("D30").Value = '714.32'
$ws.Range("D34").Value = '6.13'
$ws.Range("D35").Value = '41.35'
$ws.Range("D36").Value = '0.0₃0877'
$ws.Range("D37").Value = '61.00'
$ws.Range("D38").Value = '3.05'
$ws.Range("D39").Value = '0.0504'
$ws.Range("D40").Value = '0.147'
$ws.Range("D41").Value = '0.999'
$ws.Range("D42").Value = '0.394'
$ws.Range("D45").Value = '3.40'
$ws.Range("D48").Value = '3.42'
$ws.Range("D49").Value = '2.14'
$ws.Range("D50").Value = '145.38'
$ws.Range("D51").Value = '0.0₆0336'

# Restore default cell style now that the text values are safely stored.
$priceRange.Style = "Normal"

# --- Row 38-41 re-rank: Coin name / Link swaps ---------------------------
$ws.Range("B38").Value = 'Fetch.AI'
$ws.Range("C38").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("B41").Value = 'Dai'
$ws.Range("C41").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'

# --- Column E ("Volume(1h)") updates --------------------------------------
$ws.Range("E2").Value = '  +1.43%  '
$ws.Range("E3").Value = '  -0.88%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("E5").Value = '  -0.63%  '
$ws.Range("E6").Value = '  -1.41%  '
$ws.Range("E7").Value = '  -0.67%  '
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("E9").Value = '  +0.25%  '
$ws.Range("E10").Value = '  -0.20%  '
$ws.Range("E11").Value = '  -2.24%  '
$ws.Range("E12").Value = '  +0.00%  '
$ws.Range("E13").Value = '  +3.00%  '
$ws.Range("E14").Value = '  -0.85%  '
$ws.Range("E15").Value = '  -2.35%  '
$ws.Range("E16").Value = '  -5.48%  '
$ws.Range("E17").Value = '  -1.23%  '
$ws.Range("E18").Value = '  -0.18%  '
$ws.Range("E19").Value = '  -1.19%  '
$ws.Range("E20").Value = '  +1.35%  '
$ws.Range("E21").Value = '  -0.34%  '
$ws.Range("E22").Value = '  +4.41%  '
$ws.Range("E23").Value = '  +2.21%  '
$ws.Range("E24").Value = '  +0.62%  '
$ws.Range("E25").Value = '  +16.38%  '
$ws.Range("E26").Value = '  +11.37%  '
$ws.Range("E27").Value = '  -1.64%  '
$ws.Range("E28").Value = '  -1.71%  '
$ws.Range("E29").Value = '  +0.16%  '
$ws.Range("E30").Value = '  -2.16%  '
$ws.Range("E31").Value = '  +3.00%  '
$ws.Range("E32").Value = '  -2.12%  '
$ws.Range("E33").Value = '  +3.90%  '
$ws.Range("E34").Value = '  +12.92%  '
$ws.Range("E35").Value = '  -2.45%  '
$ws.Range("E36").Value = '  +2.67%  '
$ws.Range("E37").Value = '  +3.71%  '
$ws.Range("E38").Value = '  +19.90%  '
$ws.Range("E39").Value = '  +5.85%  '
$ws.Range("E40").Value = '  -4.66%  '
$ws.Range("E41").Value = '  +0.05%  '
$ws.Range("E42").Value = '  +15.89%  '
$ws.Range("E43").Value = '  +1.44%  '
$ws.Range("E44").Value = '  +4.98%  '
$ws.Range("E45").Value = '  +5.12%  '
$ws.Range("E46").Value = '  -1.75%  '
$ws.Range("E47").Value = '  +0.10%  '
$ws.Range("E48").Value = '  -1.46%  '
$ws.Range("E49").Value = '  -4.90%  '
$ws.Range("E50").Value = '  -1.83%  '
$ws.Range("E51").Value = '  +22.97%  '

